# Update the date heading (the single run inside the first paragraph).
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Text = "2026-01-31 Saturday"

# Update the worksheet answers, cell by cell (row, col) -> new text.
# The answer table repeats some strings verbatim in more than one cell
# (e.g. "17÷3=5, 2" appears twice, each replaced with a different new
# value, and some new values coincide with other cells' old values), so a
# document-wide Find/Replace would mismatch. Setting each table cell's own
# Range.Text targets exactly that cell while leaving its run/paragraph
# formatting (fonts, size, alignment) untouched.
$tbl = $d.Tables.Item(1)

$answers = @(
    @{Row=1;  Col=1; New="71÷2=35, 1"},
    @{Row=1;  Col=2; New="68÷8=8, 4"},
    @{Row=1;  Col=3; New="49÷7=7, 0"},
    @{Row=1;  Col=4; New="96÷4=24, 0"},
    @{Row=1;  Col=5; New="86÷2=43, 0"},

    @{Row=5;  Col=1; New="12÷3=4, 0"},
    @{Row=5;  Col=2; New="60÷9=6, 6"},
    @{Row=5;  Col=3; New="50÷9=5, 5"},
    @{Row=5;  Col=4; New="51÷7=7, 2"},
    @{Row=5;  Col=5; New="59÷6=9, 5"},

    @{Row=9;  Col=1; New="80÷2=40, 0"},
    @{Row=9;  Col=2; New="32÷5=6, 2"},
    @{Row=9;  Col=3; New="71÷9=7, 8"},
    @{Row=9;  Col=4; New="65÷7=9, 2"},
    @{Row=9;  Col=5; New="37÷4=9, 1"},

    @{Row=13; Col=1; New="68÷4=17, 0"},
    @{Row=13; Col=2; New="30÷4=7, 2"},
    @{Row=13; Col=3; New="43÷4=10, 3"},
    @{Row=13; Col=4; New="80÷9=8, 8"},
    @{Row=13; Col=5; New="21÷7=3, 0"},

    @{Row=17; Col=1; New="65÷6=10, 5"},
    @{Row=17; Col=2; New="78÷7=11, 1"},
    @{Row=17; Col=3; New="82÷9=9, 1"},
    @{Row=17; Col=4; New="17÷2=8, 1"},
    @{Row=17; Col=5; New="15÷6=2, 3"}
)

foreach ($a in $answers) {
    $cell = $tbl.Cell($a.Row, $a.Col)
    $cell.Range.Text = $a.New
}
